$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column BY (date label "21-sep"), matching style of existing
# header row (row 1).
$ws.Range("BY1").Value = "21-sep"

# New data values for column BY, rows 2-18.
$ws.Range("BY2").Value = 18.38789843037825
$ws.Range("BY3").Value = 7.7906725090427837
$ws.Range("BY4").Value = 13.385616485586558
$ws.Range("BY5").Value = 7.0382685920577623
$ws.Range("BY6").Value = 3.5133052664530298
$ws.Range("BY7").Value = 2.3634502918911746
$ws.Range("BY8").Value = 5.7246486802758279
$ws.Range("BY9").Value = 14.356981722891851
$ws.Range("BY10").Value = 15.983144021186069
$ws.Range("BY11").Value = 7.2315886699020515
$ws.Range("BY12").Value = 12.44455808387311
$ws.Range("BY13").Value = 11.563361526313006
$ws.Range("BY14").Value = 10.769519954876801
$ws.Range("BY15").Value = 6.6567925859687973
$ws.Range("BY16").Value = 11.348955959083439
$ws.Range("BY17").Value = 11.18781752150775
$ws.Range("BY18").Value = 20.676947999397147

# Update the selected cell to reflect where the user ended up after adding
# the new column (matches the saved sheet view state in the workbook).
[void]$ws.Range("BX19").Select()
